$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.909.85"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.636.27"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'211.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'23.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D12").Value = "1.867.75"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.639.14"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "'65.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "27.919.61"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'228.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'7.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").Value = "0.0₃0719"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'155.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "'6.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "1.397.34"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.76%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "'0.561"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "'0.852"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'1.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").Value = "'66.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "1.776.47"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "'88.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'7.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
